# Actualización automática 2025-06-09 15:40:08
# Insert a new client row ("IMPORTADORA ORTEGA CIA. LTDA.") right before the
# existing "INTERNEGOCIOS DE HIERRO S.A." row, in both the "VENTAS POR GRUPO"
# and "VENTA MENSUAL" sheets. All rows from that point on shift down by one,
# and the "X de 54" tally row (now the last row) is updated to "X de 55".

$wb = $excel.ActiveWorkbook

$asesor = "LINDAO ZUÑIGA BRYAN JOSE"
$nuevoCliente = "IMPORTADORA ORTEGA CIA. LTDA."

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO" (columns A:N, data rows 2..55, totals row 56)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(31).Insert()
$ws1.Range("A31").Value = $asesor
$ws1.Range("B31").Value = $nuevoCliente
$ws1.Range("C31:N31").Value = 0

$ws1.Range("C57").Value = "0 de 55"
$ws1.Range("D57").Value = "0 de 55"
$ws1.Range("E57").Value = "0 de 55"
$ws1.Range("F57").Value = "0 de 55"
$ws1.Range("G57").Value = "0 de 55"
$ws1.Range("H57").Value = "0 de 55"
$ws1.Range("I57").Value = "0 de 55"
$ws1.Range("J57").Value = "0 de 55"
$ws1.Range("K57").Value = "0 de 55"
$ws1.Range("L57").Value = "3 de 55"
$ws1.Range("M57").Value = "0 de 55"
$ws1.Range("N57").Value = "0 de 55"

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL" (columns A:G, data rows 2..55, totals row 56)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(31).Insert()
$ws2.Range("A31").Value = $asesor
$ws2.Range("B31").Value = $nuevoCliente
$ws2.Range("C31:G31").Value = 0

$wb.Save()
